$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This edit corresponds to a re-run of the KiBot/KiCad BoM generator against
# a newer KiCad version (9.0.1 -> 9.0.6+1). The underlying net-name ordering
# produced by the new version differs (same nets, different order) and a
# couple of "Net Label" picks changed because of that re-ordering.
# ---------------------------------------------------------------------------

$wsBoM = $wb.Worksheets.Item("BoM")
$wsDNF = $wb.Worksheets.Item("DNF")

# --- KiCad Version (row 6, column D) on both sheets -----------------------
$wsBoM.Range("D6").Value = "9.0.6+1"
$wsDNF.Range("D6").Value = "9.0.6+1"

# --- BoM sheet: Net Name (AC) / Net Label (AD) re-orderings ---------------

# Row 9  -> C2 C6 C7 C8 (100nF caps): GND,+5V
$wsBoM.Range("AC9").Value = "+5V,GND"
$wsBoM.Range("AD9").Value = "+5V,GND"

# Row 10 -> C5 (1uF cap): Net-(U1-UCAP),Earth
$wsBoM.Range("AC10").Value = "Earth,Net-(U1-UCAP)"
$wsBoM.Range("AD10").Value = "Earth,Net-(U1-UCAP)"

# Row 11 -> C1 (47uF cap): GND,+5V
$wsBoM.Range("AC11").Value = "+5V,GND"
$wsBoM.Range("AD11").Value = "+5V,GND"

# Row 12 -> D1 (diode): /RESET2,+5V / RESET2,+5V
$wsBoM.Range("AC12").Value = "+5V,/RESET2"
$wsBoM.Range("AD12").Value = "RESET2"

# Row 13 -> D2 D3 (LEDs): Net-(D2-A),/RXLED / RXLED
$wsBoM.Range("AC13").Value = "/RXLED,Net-(D2-A)"
$wsBoM.Range("AD13").Value = "RXLED,Net-(D2-A)"

# Row 15 -> J5/ICSP1 connector: GND,/MISO2,/RESET2,/SCK2,/MOSI2,+5V / MOSI2,+5V
$wsBoM.Range("AC15").Value = "+5V,/SCK2,/RESET2,/MISO2,GND,/MOSI2"
$wsBoM.Range("AD15").Value = "MOSI2"

# Row 16 -> J3/P1 connector: Net-(J3-Pin_*) list
$wsBoM.Range("AC16").Value = "Net-(J3-Pin_5),Net-(J3-Pin_1),Net-(J3-Pin_3),Net-(J3-Pin_2),Net-(J3-Pin_4)"
$wsBoM.Range("AD16").Value = "Net-(J3-Pin_5),Net-(J3-Pin_1),Net-(J3-Pin_3),Net-(J3-Pin_2),Net-(J3-Pin_4)"

# Row 17 -> J6/P2 connector: Net-(J6-Pin_*) list
$wsBoM.Range("AC17").Value = "Net-(J6-Pin_5),Net-(J6-Pin_3),Net-(J6-Pin_1),Net-(J6-Pin_6),Net-(J6-Pin_2),Net-(J6-Pin_4)"
$wsBoM.Range("AD17").Value = "Net-(J6-Pin_5),Net-(J6-Pin_3),Net-(J6-Pin_1),Net-(J6-Pin_6),Net-(J6-Pin_2),Net-(J6-Pin_4)"

# Row 20 -> R4: /RESET2,+5V / RESET2,+5V
$wsBoM.Range("AC20").Value = "+5V,/RESET2"
$wsBoM.Range("AD20").Value = "RESET2"

# Row 21 -> U1 (MCU): big net list + net label pick, and a shorter row height
$wsBoM.Range("AC21").Value = "Net-(U1-D-),/DTR,Earth,Net-(J4-Pin_4),Net-(J6-Pin_4),/MISO2,Net-(U1-PC0{slash}XTAL2),+5V,VBUS,Net-(J6-Pin_5),/SCK2,Net-(J4-Pin_3),Net-(J3-Pin_2),Net-(U1-XTAL1),Net-(J3-Pin_5),Net-(J6-Pin_3),Net-(J6-Pin_6),Net-(J4-Pin_1),/RESET2,Net-(J3-Pin_3),Net-(J4-Pin_2),Net-(J3-Pin_4),unconnected-(U1-PB0-Pad14),Net-(U1-D+),/RXLED,Net-(J3-Pin_1),Net-(U1-UCAP),Net-(J6-Pin_2),GND,/TXLED,/MOSI2"
$wsBoM.Range("AD21").Value = "MOSI2"
$wsBoM.Rows.Item(21).RowHeight = 105

# --- DNF sheet: Net Name (AC) / Net Label (AD) re-orderings ----------------

# Row 9  -> F1 (Polyfuse): Net-(U1-XTAL1),GND
$wsDNF.Range("AC9").Value = "GND,Net-(U1-XTAL1)"
$wsDNF.Range("AD9").Value = "GND,Net-(U1-XTAL1)"

# Row 13 -> J2 (USB_B): Net-(J2-*) list
$wsDNF.Range("AC13").Value = "Earth,Net-(J2-VBUS),Net-(J2-D-),Net-(J2-D+),Net-(J2-Shield)"
$wsDNF.Range("AD13").Value = "Earth,Net-(J2-VBUS),Net-(J2-D-),Net-(J2-D+),Net-(J2-Shield)"

# Row 14 -> R3: Net-(U1-D+),Net-(J2-D+)
$wsDNF.Range("AC14").Value = "Net-(J2-D+),Net-(U1-D+)"
$wsDNF.Range("AD14").Value = "Net-(J2-D+),Net-(U1-D+)"
